# This script updates the "cryptos" worksheet with refreshed price / volume
# data pulled from coinranking.com, matching the scheduled GitHub Actions
# data-refresh commit. For every changed cell we simply assign the new
# literal text.
#
# A number of cells in column D (Price) hold values that *look* numeric
# (e.g. "556.11", "26.24") but must remain plain text cells (t="inlineStr"
# in the original workbook) so that things like trailing zeros
# ("0.0540") and thousand-dot-separated values keep their exact textual
# form instead of being parsed/rounded as floating point numbers.
#
# To force Excel to keep such values as text we prefix them with a
# leading apostrophe (the standard "treat as text" marker) and then
# immediately reset the cell Style back to "Normal" so no left-over
# Text/quote-prefix number format lingers on the cell - the stored value
# itself no longer contains the apostrophe, only the cell's input mode
# used it.
#
# Values that are not valid numeric literals (they contain two dots like
# "62.264.78", are plain words, or are URLs) do not need this treatment
# since Excel naturally stores them as text already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').Value = '62.264.78'
$ws.Range('E2').Value = '  +2.19%  '

$ws.Range('D3').Value = '2.420.35'
$ws.Range('E3').Value = '  +2.90%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''556.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.09%  '

$ws.Range('D6').Value = '''143.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.44%  '

$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('E8').Value = '  +1.67%  '

$ws.Range('D9').Value = '2.419.75'
$ws.Range('E9').Value = '  +2.83%  '

$ws.Range('E10').Value = '  +3.99%  '

$ws.Range('E11').Value = '  -0.62%  '

$ws.Range('E12').Value = '  +1.35%  '

$ws.Range('E13').Value = '  +2.10%  '

$ws.Range('D14').Value = '''26.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.06%  '

$ws.Range('E15').Value = '  +8.20%  '

$ws.Range('D16').Value = '2.858.67'
$ws.Range('E16').Value = '  +2.91%  '

$ws.Range('D17').Value = '62.249.63'
$ws.Range('E17').Value = '  +2.25%  '

$ws.Range('D18').Value = '2.420.04'
$ws.Range('E18').Value = '  +5.04%  '

$ws.Range('E19').Value = '  +3.99%  '

$ws.Range('E20').Value = '  +2.08%  '

$ws.Range('D21').Value = '''323.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.21%  '

$ws.Range('E22').Value = '  +2.14%  '

$ws.Range('E23').Value = '  +0.27%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''64.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.34%  '

$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').Value = '''1.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.41%  '

$ws.Range('D26').Value = '''9.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.01%  '

$ws.Range('D27').Value = '''578.69'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.07%  '

$ws.Range('D28').Value = '2.541.64'
$ws.Range('E28').Value = '  +2.91%  '

$ws.Range('E29').Value = '  +0.09%  '

$ws.Range('E30').Value = '  +4.66%  '

$ws.Range('E31').Value = '  +8.11%  '

$ws.Range('E32').Value = '  +5.61%  '

$ws.Range('E33').Value = '  +1.45%  '

$ws.Range('E34').Value = '  +4.12%  '

$ws.Range('E35').Value = '  +3.30%  '

$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '''0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '''5.67'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.96%  '

$ws.Range('E38').Value = '  +3.97%  '

$ws.Range('D39').Value = '''0.384'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.12%  '

$ws.Range('E40').Value = '  +3.30%  '

$ws.Range('D41').Value = '''18.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.97%  '

$ws.Range('D42').Value = '''148.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.69%  '

$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('D44').Value = '''41.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.66%  '

$ws.Range('E45').Value = '  +13.15%  '

$ws.Range('D46').Value = '''150.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.50%  '

$ws.Range('D47').Value = '''3.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.68%  '

$ws.Range('D48').Value = '''0.0540'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.76%  '

$ws.Range('D49').Value = '''20.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.23%  '

$ws.Range('E50').Value = '  +3.40%  '

$ws.Range('D51').Value = '''0.0918'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.78%  '
